# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.050.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.58"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.18"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.64%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.113.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.845.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.678"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.70"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.058.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0794"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("E29").Value = "  +12.62%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0557"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("B31").Value = "BinanceUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  +24.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.764"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("E38").Value = "  +12.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "92.18"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0202"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.351.55"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.67"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +86.92%  "
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.029.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +17.47%  "
$ws.Range("E51").Value = "  +0.68%  "
